# Rebuild Sheet1 as a data-driven login credentials table (email/password),
# replacing the old Nico/Ardy/Jakarta... sample data, and mark the email
# column with Hyperlink-styled cells (some of them real mailto hyperlinks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so stale shared strings are dropped.
$ws.Cells.Clear()

# --- Body rows first (this reproduces the shared-string insertion order
#     seen in the target file: email variants, then admin variants, with
#     the "email"/"password" header strings added last). ---
$ws.Range("A2").Value = "Admin@yourstore.com"
$ws.Range("A3").Value = "ADMIN@YOURSTORE.COM"
$ws.Range("A4").Value = "ADMIN@yourstore.com"
$ws.Range("A5").Value = "admin@yourstore.com"
$ws.Range("B2").Value = "admin"
$ws.Range("B3").Value = "admin"
$ws.Range("B4").Value = "admin"
$ws.Range("B5").Value = "admin"

$ws.Range("A6").Value = "ADMIN@YOURSTORE.COM"
$ws.Range("B6").Value = "admin1"
$ws.Range("A7").Value = "ADMIN@yourstore.com"
$ws.Range("B7").Value = "admin2"
$ws.Range("A8").Value = "Admin@yourstore.com"
$ws.Range("B8").Value = "admin3"
$ws.Range("A9").Value = "admin@yourstore.com"
$ws.Range("B9").Value = "admin4"

# --- Header row last ---
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "password"

# --- Real hyperlinks on the first occurrence of each distinct e-mail ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:admin@yourstore.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:admin@yourstore.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:admin@yourstore.com")
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:admin@yourstore.com")
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:admin@yourstore.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:admin@yourstore.com")

# --- Rows 6 & 7 keep the Hyperlink look (underline + theme color) without
#     being live links. ---
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("A7").Style = "Hyperlink"

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 23

# --- Final selection lands on B1 ---
$ws.Range("B1").Select() | Out-Null
